$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI LR-pair table is expanded from a single Sending->Target cluster
# combination (FAPs -> M2) into all four Sending/Target cluster combinations
# (FAPs/Neutro x M1/M2) for the Ccl20/Ccr6 ligand-receptor pair, per Dr Hou's
# advice. Values are written column-by-column so that brand new category
# labels ("Neutro", "M1") are interned in the same relative order as the
# original authoring tool used.

# Column A: Sending cluster
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "Neutro"
$ws.Range("A5").Value = "Neutro"

# Column B: Ligand symbol
$ws.Range("B2").Value = "Ccl20"
$ws.Range("B3").Value = "Ccl20"
$ws.Range("B4").Value = "Ccl20"
$ws.Range("B5").Value = "Ccl20"

# Column C: Receptor symbol
$ws.Range("C2").Value = "Ccr6"
$ws.Range("C3").Value = "Ccr6"
$ws.Range("C4").Value = "Ccr6"
$ws.Range("C5").Value = "Ccr6"

# Column D: Target cluster
$ws.Range("D2").Value = "M1"
$ws.Range("D3").Value = "M2"
$ws.Range("D4").Value = "M1"
$ws.Range("D5").Value = "M2"

# Row 2: FAPs -> Ccl20/Ccr6 -> M1
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.019213
$ws.Range("H2").Value = 3.057639
$ws.Range("I2").Value = 0.928181872270981
$ws.Range("J2").Value = 0.928181872270981
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04774999999999999
$ws.Range("N2").Value = 0.14325
$ws.Range("O2").Value = 0.0255348367365819
$ws.Range("P2").Value = 0.0255348367365819
$ws.Range("Q2").Value = 0.04866742074999999
$ws.Range("R2").Value = 0.43800678675
$ws.Range("S2").Value = 0.02370097257029442
$ws.Range("T2").Value = 0.02370097257029442

# Row 3: FAPs -> Ccl20/Ccr6 -> M2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.019213
$ws.Range("H3").Value = 3.057639
$ws.Range("I3").Value = 0.928181872270981
$ws.Range("J3").Value = 0.928181872270981
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.822244333333334
$ws.Range("N3").Value = 5.466733000000001
$ws.Range("O3").Value = 0.9744651632634181
$ws.Range("P3").Value = 0.974465163263418
$ws.Range("Q3").Value = 1.857255113709667
$ws.Range("R3").Value = 16.715296023387
$ws.Range("S3").Value = 0.9044808997006867
$ws.Range("T3").Value = 0.9044808997006866

# Row 4: Neutro -> Ccl20/Ccr6 -> M1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07886166666666666
$ws.Range("H4").Value = 0.236585
$ws.Range("I4").Value = 0.07181812772901904
$ws.Range("J4").Value = 0.07181812772901904
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04774999999999999
$ws.Range("N4").Value = 0.14325
$ws.Range("O4").Value = 0.0255348367365819
$ws.Range("P4").Value = 0.0255348367365819
$ws.Range("Q4").Value = 0.003765644583333333
$ws.Range("R4").Value = 0.03389080125
$ws.Range("S4").Value = 0.001833864166287487
$ws.Range("T4").Value = 0.001833864166287487

# Row 5: Neutro -> Ccl20/Ccr6 -> M2
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07886166666666666
$ws.Range("H5").Value = 0.236585
$ws.Range("I5").Value = 0.07181812772901904
$ws.Range("J5").Value = 0.07181812772901904
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.822244333333334
$ws.Range("N5").Value = 5.466733000000001
$ws.Range("O5").Value = 0.9744651632634181
$ws.Range("P5").Value = 0.974465163263418
$ws.Range("Q5").Value = 0.1437052252005556
$ws.Range("R5").Value = 1.293347026805
$ws.Range("S5").Value = 0.06998426356273155
$ws.Range("T5").Value = 0.06998426356273155
